$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values (A:N) for every row involved in the reshuffle,
# before any writes happen, since several rows swap/rotate with each other.
# Value2 is used because multi-cell Value reads are not materialized reliably
# in this host; Value2 returns a proper 2-D array.
$row2 = $ws.Range("A2:N2").Value2
$row3 = $ws.Range("A3:N3").Value2
$row12 = $ws.Range("A12:N12").Value2
$row13 = $ws.Range("A13:N13").Value2
$row16 = $ws.Range("A16:N16").Value2
$row17 = $ws.Range("A17:N17").Value2
$row18 = $ws.Range("A18:N18").Value2
$row19 = $ws.Range("A19:N19").Value2
$row40 = $ws.Range("A40:N40").Value2
$row41 = $ws.Range("A41:N41").Value2
$row42 = $ws.Range("A42:N42").Value2
$row43 = $ws.Range("A43:N43").Value2
$row44 = $ws.Range("A44:N44").Value2
$row45 = $ws.Range("A45:N45").Value2
$row46 = $ws.Range("A46:N46").Value2
$row47 = $ws.Range("A47:N47").Value2
$row52 = $ws.Range("A52:N52").Value2
$row53 = $ws.Range("A53:N53").Value2
$row58 = $ws.Range("A58:N58").Value2
$row60 = $ws.Range("A60:N60").Value2
$row61 = $ws.Range("A61:N61").Value2

# Write each captured row into its new destination per the commit's reordering
$ws.Range("A2:N2").Value2 = $row3
$ws.Range("A3:N3").Value2 = $row2
$ws.Range("A12:N12").Value2 = $row13
$ws.Range("A13:N13").Value2 = $row12
$ws.Range("A16:N16").Value2 = $row17
$ws.Range("A17:N17").Value2 = $row16
$ws.Range("A18:N18").Value2 = $row19
$ws.Range("A19:N19").Value2 = $row18
$ws.Range("A40:N40").Value2 = $row47
$ws.Range("A41:N41").Value2 = $row44
$ws.Range("A42:N42").Value2 = $row43
$ws.Range("A43:N43").Value2 = $row45
$ws.Range("A44:N44").Value2 = $row42
$ws.Range("A45:N45").Value2 = $row40
$ws.Range("A46:N46").Value2 = $row41
$ws.Range("A47:N47").Value2 = $row46
$ws.Range("A52:N52").Value2 = $row53
$ws.Range("A53:N53").Value2 = $row52
$ws.Range("A58:N58").Value2 = $row60
$ws.Range("A60:N60").Value2 = $row61
$ws.Range("A61:N61").Value2 = $row58
